$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Add the (hidden) DropdownOptions sheet after Sheet1 ---
$ddSheet = $wb.Worksheets.Add($null, $ws1)
$ddSheet.Name = "DropdownOptions"

$options = @(
    "0% - 10%: Foundation completed: Groundwork finished; no vertical structure yet.",
    "11% - 25%: Structure and rough-in started: Structural framing in progress; initial MEP rough-in.",
    "26% - 50%: Structure erected, partial roofing: Building shape defined; roof and systems advancing.",
    "51% - 75%: Exterior sealed, interior work underway: Enclosed structure; painting, flooring, and testing begin.",
    "76% - 90%: Final finishes and inspections: Systems tested; ",
    "91% - 99%: Final touches and punch list: Minor adjustments; final inspections and approvals.",
    "100% - Construction complete: Ready for handover and occupancy."
)

for ($i = 0; $i -lt $options.Count; $i++) {
    $ddSheet.Cells.Item($i + 1, 1).Value = $options[$i]
}

# --- New "Status" header column on Sheet1 ---
$ws1.Range("AH1").Value = "Status"

# --- Remove the stray empty cells that used to sit at AE2 / AF2 ---
$ws1.Range("AE2").ClearContents()
$ws1.Range("AF2").ClearContents()

# --- Drop-down list validation on AH2, sourced from DropdownOptions ---
$validation = $ws1.Range("AH2").Validation
$validation.Add(3, 1, 1, '=DropdownOptions!$A$1:$A$7')
$validation.ShowInput = $false
$validation.ShowError = $false

# --- Hide the helper sheet ---
$ddSheet.Visible = $false

# Keep the original sheet active/selected
$ws1.Activate()
